$wb = $excel.ActiveWorkbook
$gens = $wb.Worksheets.Item("gens")
$lines = $wb.Worksheets.Item("lines")
$bus = $wb.Worksheets.Item("bus")

# --- gens sheet: production column (C) updates ---
$gens.Cells.Item(10, 3).Value = 100
$gens.Cells.Item(12, 3).Value = 100
$gens.Cells.Item(24, 3).Value = 400

# --- lines sheet: flow_loads / shadow price updates ---
$lines.Cells.Item(2, 4).Value = 0.12152197
$lines.Cells.Item(3, 4).Value = -0.55512114
$lines.Cells.Item(4, 4).Value = -0.51418495
$lines.Cells.Item(5, 4).Value = -0.4445064
$lines.Cells.Item(6, 4).Value = 2.3385967
$lines.Cells.Item(7, 4).Value = -0.13186853
$lines.Cells.Item(8, 4).Value = -0.12791242
$lines.Cells.Item(9, 4).Value = -0.36400524
$lines.Cells.Item(10, 4).Value = -0.53233265
$lines.Cells.Item(11, 4).Value = 0.74299165
$lines.Cells.Item(11, 5).Value = 4.2496279
$lines.Cells.Item(11, 6).Value = 0
$lines.Cells.Item(12, 3).Value = 157.5
$lines.Cells.Item(12, 6).Value = -0.4092255
$lines.Cells.Item(13, 3).Value = -23.805448
$lines.Cells.Item(13, 4).Value = 0.17976396
$lines.Cells.Item(14, 3).Value = 10.305448
$lines.Cells.Item(14, 4).Value = -0.17976396
$lines.Cells.Item(15, 3).Value = -164.10417
$lines.Cells.Item(15, 4).Value = -0.14416418
$lines.Cells.Item(16, 3).Value = -156.34132
$lines.Cells.Item(16, 4).Value = -0.15140753
$lines.Cells.Item(17, 3).Value = -231.10771
$lines.Cells.Item(17, 4).Value = 0.21536374
$lines.Cells.Item(18, 3).Value = -223.34487
$lines.Cells.Item(18, 4).Value = 0.20812039
$lines.Cells.Item(19, 3).Value = -304.34817
$lines.Cells.Item(19, 4).Value = 0.0045863754
$lines.Cells.Item(20, 3).Value = -90.863709
$lines.Cells.Item(20, 4).Value = 0.031586706
$lines.Cells.Item(21, 3).Value = -317.93315
$lines.Cells.Item(21, 4).Value = 0.011829732
$lines.Cells.Item(22, 3).Value = -61.753037
$lines.Cells.Item(22, 4).Value = 0.041583927
$lines.Cells.Item(23, 3).Value = 106.56031
$lines.Cells.Item(23, 4).Value = 0.029754195
$lines.Cells.Item(24, 3).Value = -284.86371
$lines.Cells.Item(24, 4).Value = 0.044371801
$lines.Cells.Item(25, 3).Value = 192.90573
$lines.Cells.Item(25, 4).Value = 0.0070221036
$lines.Cells.Item(26, 3).Value = -314.75383
$lines.Cells.Item(26, 4).Value = -0.047427881
$lines.Cells.Item(27, 3).Value = -314.75383
$lines.Cells.Item(27, 4).Value = -0.047427881
$lines.Cells.Item(28, 4).Value = 0.079183882
$lines.Cells.Item(29, 3).Value = -450
$lines.Cells.Item(29, 4).Value = 0.050331629
$lines.Cells.Item(29, 5).Value = 0.14999297
$lines.Cells.Item(30, 3).Value = 413.04202
$lines.Cells.Item(30, 4).Value = -0.017726158
$lines.Cells.Item(31, 3).Value = -291.605
$lines.Cells.Item(31, 4).Value = 0.023442924
$lines.Cells.Item(32, 3).Value = -158.395
$lines.Cells.Item(32, 4).Value = 0.027440417
$lines.Cells.Item(33, 3).Value = -112.3025
$lines.Cells.Item(33, 4).Value = 0.021768429
$lines.Cells.Item(34, 3).Value = -112.3025
$lines.Cells.Item(34, 4).Value = 0.021768429
$lines.Cells.Item(35, 3).Value = 116.02101
$lines.Cells.Item(35, 4).Value = -0.015414051
$lines.Cells.Item(36, 3).Value = 116.02101
$lines.Cells.Item(36, 4).Value = -0.015414051
$lines.Cells.Item(37, 3).Value = 52.021009
$lines.Cells.Item(37, 4).Value = -0.008477727799999999
$lines.Cells.Item(38, 3).Value = 52.021009
$lines.Cells.Item(38, 4).Value = -0.008477727799999999
$lines.Cells.Item(39, 3).Value = -141.605
$lines.Cells.Item(39, 4).Value = -0.017770937

# --- bus sheet: bus_lmp / node_theta updates ---
$bus.Cells.Item(2, 2).Value = 110.56865
$bus.Cells.Item(2, 3).Value = 0
$bus.Cells.Item(3, 2).Value = 110.69017
$bus.Cells.Item(3, 3).Value = -0.0034986073
$bus.Cells.Item(4, 2).Value = 110.01353
$bus.Cells.Item(4, 3).Value = 0.19355457
$bus.Cells.Item(5, 2).Value = 110.24567
$bus.Cells.Item(5, 3).Value = 0.060649027
$bus.Cells.Item(6, 2).Value = 110.05447
$bus.Cells.Item(6, 3).Value = 0.035069326
$bus.Cells.Item(7, 2).Value = 113.02877
$bus.Cells.Item(7, 3).Value = 0.037781393
$bus.Cells.Item(8, 2).Value = 109.29267
$bus.Cells.Item(8, 3).Value = 0.24693538
$bus.Cells.Item(9, 2).Value = 109.7019
$bus.Cells.Item(9, 3).Value = 0.15086038
$bus.Cells.Item(10, 2).Value = 109.88166
$bus.Cells.Item(10, 3).Value = 0.19013937
$bus.Cells.Item(11, 2).Value = 109.52213
$bus.Cells.Item(11, 3).Value = 0.13385639
$bus.Cells.Item(12, 2).Value = 109.7375
$bus.Cells.Item(12, 3).Value = 0.32798687
$bus.Cells.Item(13, 2).Value = 109.73025
$bus.Cells.Item(13, 3).Value = 0.32146608
$bus.Cells.Item(14, 2).Value = 109.74208
$bus.Cells.Item(14, 3).Value = 0.47407399
$bus.Cells.Item(15, 2).Value = 109.76908
$bus.Cells.Item(15, 3).Value = 0.36614963
$bus.Cells.Item(16, 2).Value = 109.80643
$bus.Cells.Item(16, 3).Value = 0.5670131899999999
$bus.Cells.Item(17, 2).Value = 109.81346
$bus.Cells.Item(17, 3).Value = 0.53421922
$bus.Cells.Item(18, 2).Value = 109.71379
$bus.Cells.Item(18, 3).Value = 0.65121922
$bus.Cells.Item(19, 2).Value = 109.73724
$bus.Cells.Item(19, 3).Value = 0.69204392
$bus.Cells.Item(20, 2).Value = 109.79573
$bus.Cells.Item(20, 3).Value = 0.43921955
$bus.Cells.Item(21, 2).Value = 109.78032
$bus.Cells.Item(21, 3).Value = 0.39281115
$bus.Cells.Item(22, 2).Value = 109.75901
$bus.Cells.Item(22, 3).Value = 0.72124257
$bus.Cells.Item(23, 2).Value = 109.74124
$bus.Cells.Item(23, 3).Value = 0.8175339700000001
$bus.Cells.Item(24, 2).Value = 109.77184
$bus.Cells.Item(24, 3).Value = 0.38136653
$bus.Cells.Item(25, 2).Value = 109.88562
$bus.Cells.Item(25, 3).Value = 0.42422019

